$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '23.966.92' },
    @{ Cell = 'E2'; Value = '  -0.74%  ' },
    @{ Cell = 'D3'; Value = '1.651.34' },
    @{ Cell = 'E3'; Value = '  +0.71%  ' },
    @{ Cell = 'D4'; Value = '1.000' },
    @{ Cell = 'E4'; Value = '  +0.10%  ' },
    @{ Cell = 'D5'; Value = '308.89' },
    @{ Cell = 'E5'; Value = '  -0.81%  ' },
    @{ Cell = 'D6'; Value = '1.001' },
    @{ Cell = 'E6'; Value = '  +0.16%  ' },
    @{ Cell = 'D7'; Value = '0.3885' },
    @{ Cell = 'E7'; Value = '  -1.30%  ' },
    @{ Cell = 'D8'; Value = '0.3824' },
    @{ Cell = 'E8'; Value = '  -1.72%  ' },
    @{ Cell = 'D9'; Value = '51.80' },
    @{ Cell = 'E9'; Value = '  +0.84%  ' },
    @{ Cell = 'D10'; Value = '1.350' },
    @{ Cell = 'E10'; Value = '  -2.86%  ' },
    @{ Cell = 'D11'; Value = '1.000' },
    @{ Cell = 'E11'; Value = '  +0.12%  ' },
    @{ Cell = 'D12'; Value = '0.08453' },
    @{ Cell = 'E12'; Value = '  -1.15%  ' },
    @{ Cell = 'D13'; Value = '23.86' },
    @{ Cell = 'E13'; Value = '  -0.93%  ' },
    @{ Cell = 'D14'; Value = '7.082' },
    @{ Cell = 'E14'; Value = '  -1.77%  ' },
    @{ Cell = 'D15'; Value = '7.925' },
    @{ Cell = 'E15'; Value = '  +3.10%  ' },
    @{ Cell = 'D16'; Value = '0.00001313' },
    @{ Cell = 'E16'; Value = '  -2.33%  ' },
    @{ Cell = 'D17'; Value = '1.645.83' },
    @{ Cell = 'E17'; Value = '  +0.06%  ' },
    @{ Cell = 'D18'; Value = '94.76' },
    @{ Cell = 'E18'; Value = '  -0.41%  ' },
    @{ Cell = 'D19'; Value = '0.06955' },
    @{ Cell = 'E19'; Value = '  +0.36%  ' },
    @{ Cell = 'D20'; Value = '19.69' },
    @{ Cell = 'E20'; Value = '  -3.22%  ' },
    @{ Cell = 'D21'; Value = '6.972' },
    @{ Cell = 'E21'; Value = '  +0.59%  ' },
    @{ Cell = 'E22'; Value = '  +0.24%  ' },
    @{ Cell = 'D23'; Value = '13.78' },
    @{ Cell = 'E23'; Value = '  +1.49%  ' },
    @{ Cell = 'D24'; Value = '23.961.90' },
    @{ Cell = 'E24'; Value = '  -0.74%  ' },
    @{ Cell = 'D25'; Value = '2.449' },
    @{ Cell = 'E25'; Value = '  -0.64%  ' },
    @{ Cell = 'D26'; Value = '3.017' },
    @{ Cell = 'E26'; Value = '  +3.65%  ' },
    @{ Cell = 'D27'; Value = '22.15' },
    @{ Cell = 'E27'; Value = '  -1.11%  ' },
    @{ Cell = 'D28'; Value = '151.78' },
    @{ Cell = 'E28'; Value = '  -3.99%  ' },
    @{ Cell = 'D29'; Value = '5.391' },
    @{ Cell = 'E29'; Value = '  +0.41%  ' },
    @{ Cell = 'D30'; Value = '138.57' },
    @{ Cell = 'E30'; Value = '  -1.86%  ' },
    @{ Cell = 'D31'; Value = '7.803' },
    @{ Cell = 'E31'; Value = '  -2.73%  ' },
    @{ Cell = 'D32'; Value = '2.532' },
    @{ Cell = 'E32'; Value = '  +0.30%  ' },
    @{ Cell = 'D33'; Value = '1.828.57' },
    @{ Cell = 'E33'; Value = '  +0.73%  ' },
    @{ Cell = 'D34'; Value = '1.040' },
    @{ Cell = 'E34'; Value = '  +2.68%  ' },
    @{ Cell = 'D35'; Value = '0.08067' },
    @{ Cell = 'E35'; Value = '  -1.77%  ' },
    @{ Cell = 'D36'; Value = '0.02945' },
    @{ Cell = 'E36'; Value = '  +0.66%  ' },
    @{ Cell = 'D37'; Value = '6.648' },
    @{ Cell = 'E37'; Value = '  -1.88%  ' },
    @{ Cell = 'D38'; Value = '10.79' },
    @{ Cell = 'E38'; Value = '  +1.65%  ' },
    @{ Cell = 'D39'; Value = '0.2667' },
    @{ Cell = 'E39'; Value = '  -1.46%  ' },
    @{ Cell = 'D40'; Value = '0.09100' },
    @{ Cell = 'E40'; Value = '  -1.56%  ' },
    @{ Cell = 'D41'; Value = '0.7586' },
    @{ Cell = 'E41'; Value = '  -1.34%  ' },
    @{ Cell = 'D42'; Value = '13.40' },
    @{ Cell = 'E42'; Value = '  -4.82%  ' },
    @{ Cell = 'D43'; Value = '1.421' },
    @{ Cell = 'E43'; Value = '  -0.70%  ' },
    @{ Cell = 'D44'; Value = '16.27' },
    @{ Cell = 'E44'; Value = '  -0.53%  ' },
    @{ Cell = 'D45'; Value = '0.7009' },
    @{ Cell = 'E45'; Value = '  -0.44%  ' },
    @{ Cell = 'D46'; Value = '2.463' },
    @{ Cell = 'E46'; Value = '  -1.81%  ' },
    @{ Cell = 'D47'; Value = '4.079' },
    @{ Cell = 'E47'; Value = '  -0.80%  ' },
    @{ Cell = 'D48'; Value = '1.001' },
    @{ Cell = 'E48'; Value = '  +0.17%  ' },
    @{ Cell = 'D49'; Value = '0.08308' },
    @{ Cell = 'E49'; Value = '  -0.52%  ' },
    @{ Cell = 'D50'; Value = '135.06' },
    @{ Cell = 'E50'; Value = '  -0.85%  ' },
    @{ Cell = 'D51'; Value = '1.209' },
    @{ Cell = 'E51'; Value = '  -2.44%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = "'" + $u.Value
    $ws.Range($u.Cell).ClearFormats()
}
